$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.272.92"
$ws.Range("E2").Value = "  +2.82%  "
$ws.Range("D3").Value = "2.596.98"
$ws.Range("E3").Value = "  +1.71%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "535.56"
$ws.Range("E5").Value = "  +4.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.96"
$ws.Range("E6").Value = "  +2.66%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("E8").Value = "  +2.03%  "
$ws.Range("D9").Value = "2.612.14"
$ws.Range("E9").Value = "  +2.18%  "
$ws.Range("E10").Value = "  +0.61%  "
$ws.Range("E11").Value = "  +4.40%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.336"
$ws.Range("E12").Value = "  +4.15%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.135"
$ws.Range("E13").Value = "  +2.71%  "
$ws.Range("D14").Value = "3.057.42"
$ws.Range("E14").Value = "  +1.71%  "
$ws.Range("D15").Value = "59.204.96"
$ws.Range("E15").Value = "  +2.61%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.61"
$ws.Range("E16").Value = "  +3.54%  "
$ws.Range("D17").Value = "2.649.20"
$ws.Range("E17").Value = "  +1.72%  "
$ws.Range("E18").Value = "  +3.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "346.30"
$ws.Range("E19").Value = "  +4.79%  "
$ws.Range("E20").Value = "  +2.14%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.12"
$ws.Range("E21").Value = "  +1.49%  "
$ws.Range("E22").Value = "  +1.41%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "67.11"
$ws.Range("E24").Value = "  +2.70%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.169"
$ws.Range("E25").Value = "  +3.42%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.408"
$ws.Range("E26").Value = "  +3.24%  "
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.23"
$ws.Range("E28").Value = "  +5.48%  "
$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").Value = "0.0₃0740"
$ws.Range("E29").Value = "  +6.12%  "
$ws.Range("B30").Value = "USDe"
$ws.Range("C30").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("E31").Value = "  +5.67%  "
$ws.Range("E32").Value = "  -1.16%  "
$ws.Range("E33").Value = "  +1.97%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "149.57"
$ws.Range("E34").Value = "  +0.55%  "
$ws.Range("E35").Value = "  +3.81%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.13"
$ws.Range("E36").Value = "  +2.58%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "36.94"
$ws.Range("E37").Value = "  +2.94%  "
$ws.Range("E38").Value = "  +5.83%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.845"
$ws.Range("E39").Value = "  +3.34%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.833"
$ws.Range("E40").Value = "  +1.72%  "
$ws.Range("E41").Value = "  +2.25%  "
$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "278.78"
$ws.Range("E42").Value = "  +3.59%  "
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.997"
$ws.Range("E43").Value = "  +0.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.599"
$ws.Range("E44").Value = "  +3.20%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.76"
$ws.Range("E45").Value = "  +0.63%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0963"
$ws.Range("E46").Value = "  +3.03%  "
$ws.Range("E47").Value = "  +2.48%  "
$ws.Range("D48").Value = "1.943.10"
$ws.Range("E48").Value = "  -0.88%  "
$ws.Range("E49").Value = "  +3.57%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.38"
$ws.Range("E50").Value = "  +4.89%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.50"
$ws.Range("E51").Value = "  +3.94%  "
